$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Page3")

# Widen column I (9th column) so the stored width goes from 14 to 15
# (ColumnWidth uses a pixel-rounded character scale, 14.15 serializes to width="15")
$ws.Columns.Item(9).ColumnWidth = 14.15

# Update certification paragraph text
$ws.Range("A18").Value = "This is to certify that Ms. IMEE JANINE O. ABALON has graduated with the degree of"
$ws.Range("A19").Value = "Bachelor of Science in Industrial Education (BSIE), major in Food and Service Management"
$ws.Range("A20").Value = "from College Of Industrial Technology, East Campus, Legazpi City, Albay on March 28, 2015"
$ws.Range("A21").Value = "per Board of Regents Referendum No. 02-A, s. 2015 having a General Weighted Average (GWA)"
$ws.Range("A24").Value = "Issued this 15th day of July, 2021 upon the request of interested party for reference purposes."

# Row 25 text moved into row 24; unmerge A25:I25 and clear its content,
# matching the plain (unindented) style used by the rest of that row
$ws.Range("A25:I25").UnMerge()
$ws.Range("A25").Value = $null
$ws.Range("B25").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update signatory name
$ws.Range("G30").Value = "CORAZON N. BAZAR"
